$d = $word.ActiveDocument

# 1. Shorten the title text.
$d.Content.Find.Execute(
    "Gravity Optimizer: A Mechanical View on Optimization in Deep Learning",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Gravity Optimizer", 2
)

# 2. Word re-drops the "_GoBack" bookmark (tracking the last edit location)
#    onto the empty paragraph right after the title -- this also moves it
#    away from its old spot near the end of the document, and bumps the
#    "_3.2_Datasets" bookmark's id out of the way.
$lastEditPara = $d.Paragraphs.Item(2)
$d.Bookmarks.Add("_GoBack", $lastEditPara.Range)
